$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

$ws.Cells.Item($row, 1).Value = "2024-10-01 00:00:00"
$ws.Cells.Item($row, 2).Value = 75650
$ws.Cells.Item($row, 3).Value = 10756.89
$ws.Cells.Item($row, 4).Value = 9519.370000000001
$ws.Cells.Item($row, 5).Value = 7.0284
